$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.938.68'
$ws.Range('E2').Value = '  +1.97%  '

$ws.Range('D3').Value = '3.375.93'
$ws.Range('E3').Value = '  +7.39%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '628.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.53%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.39'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +23.85%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.392'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.73%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.863'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.52%  '

$ws.Range('B11').Value = 'LidoStakedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D11').Value = '3.373.72'
$ws.Range('E11').Value = '  +7.40%  '

$ws.Range('E12').Value = '  +0.20%  '

$ws.Range('D13').Value = '98.768.06'
$ws.Range('E13').Value = '  +2.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.11'
$ws.Range('D14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000248'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.00%  '

$ws.Range('D16').Value = '3.961.85'
$ws.Range('E16').Value = '  +6.28%  '

$ws.Range('E17').Value = '  -0.54%  '

$ws.Range('D18').Value = '3.384.21'
$ws.Range('E18').Value = '  +8.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.11%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '490.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.90%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000210'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.41%  '

$ws.Range('E24').Value = '  +6.50%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.48%  '

$ws.Range('D28').Value = '3.558.14'
$ws.Range('E28').Value = '  +7.58%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.280'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.54%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.15%  '

$ws.Range('E31').Value = '  +8.88%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.131'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.33%  '

$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +11.41%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.64%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '28.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.35%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.32'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.93%  '

$ws.Range('E37').Value = '  -1.85%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.96'
$ws.Range('D38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '499.68'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.13%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.460'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '24.90'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.80%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.77'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.26'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.52%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.789'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +11.73%  '

$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '160.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.54%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.843'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +13.32%  '

$ws.Range('E50').Value = '  +2.92%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.81%  '
